# Apply edit: refresh sensor reading rows 2-5 with a new data sample (34 junction
# columns), drop the now-obsolete 6th data row, and resize several columns that
# needed extra room ("custom accuracy" + swapped-in 1000-row dataset sample).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Overwrite data rows 2-5 (columns A..AH) with the new reading values ---
$row2 = New-Object 'object[,]' 1,34
$row2[0,0] = 45042.50694444445
$row2[0,1] = 18.737
$row2[0,2] = 12.527
$row2[0,3] = 4.029
$row2[0,4] = 39.827
$row2[0,5] = 31.89
$row2[0,6] = 14.745
$row2[0,7] = 46.499
$row2[0,8] = 22.687
$row2[0,9] = 9.44
$row2[0,10] = 14.208
$row2[0,11] = 15.679
$row2[0,12] = 16.314
$row2[0,13] = 4.706
$row2[0,14] = 14.663
$row2[0,15] = 20.445
$row2[0,16] = 12.568
$row2[0,17] = 3.441
$row2[0,18] = 2.244
$row2[0,19] = 215.823
$row2[0,20] = 40.754
$row2[0,21] = 13.534
$row2[0,22] = 26.795
$row2[0,23] = 13.632
$row2[0,24] = 3.069
$row2[0,25] = 23.634
$row2[0,26] = 11.955
$row2[0,27] = 10.86
$row2[0,28] = 12.744
$row2[0,29] = 16.181
$row2[0,30] = 3.455
$row2[0,31] = 41.233
$row2[0,32] = 7.422
$row2[0,33] = 16.92
$ws.Range("A2:AH2").Value2 = $row2

$row3 = New-Object 'object[,]' 1,34
$row3[0,0] = 45042.51388888889
$row3[0,1] = 5.765
$row3[0,2] = 3.615
$row3[0,3] = 1.454
$row3[0,4] = 12.323
$row3[0,5] = 9.539
$row3[0,6] = 4.538
$row3[0,7] = 21.211
$row3[0,8] = 6.981
$row3[0,9] = 2.836
$row3[0,10] = 4.019
$row3[0,11] = 4.9
$row3[0,12] = 5.003
$row3[0,13] = 1.456
$row3[0,14] = 4.512
$row3[0,15] = 6.245
$row3[0,16] = 4.159
$row3[0,17] = 1.459
$row3[0,18] = 0.755
$row3[0,19] = 61.386
$row3[0,20] = 12.859
$row3[0,21] = 4.164
$row3[0,22] = 8.242000000000001
$row3[0,23] = 4.129
$row3[0,24] = 1.195
$row3[0,25] = 9.949999999999999
$row3[0,26] = 3.678
$row3[0,27] = 3.508
$row3[0,28] = 4.086
$row3[0,29] = 4.974
$row3[0,30] = 1.265
$row3[0,31] = 19.661
$row3[0,32] = 2.17
$row3[0,33] = 5.208
$ws.Range("A3:AH3").Value2 = $row3

$row4 = New-Object 'object[,]' 1,34
$row4[0,0] = 45042.52083333334
$row4[0,1] = 0.373
$row4[0,2] = 0.181
$row4[0,3] = 0.787
$row4[0,4] = 0.889
$row4[0,5] = 0
$row4[0,6] = 0
$row4[0,7] = 5.975
$row4[0,8] = 0.582
$row4[0,9] = 0.107
$row4[0,10] = 0.16
$row4[0,11] = 0.261
$row4[0,12] = 0
$row4[0,13] = 0
$row4[0,14] = 0.376
$row4[0,15] = 0.531
$row4[0,16] = 0.587
$row4[0,17] = 0.912
$row4[0,18] = 0.32
$row4[0,19] = 0
$row4[0,20] = 1.461
$row4[0,21] = 0.347
$row4[0,22] = 0.827
$row4[0,23] = 0.286
$row4[0,24] = 0.461
$row4[0,25] = 2.541
$row4[0,26] = 0.307
$row4[0,27] = 0.399
$row4[0,28] = 0.429
$row4[0,29] = 0.284
$row4[0,30] = 0.772
$row4[0,31] = 5.875
$row4[0,32] = 0.041
$row4[0,33] = 0.454
$ws.Range("A4:AH4").Value2 = $row4

$row5 = New-Object 'object[,]' 1,34
$row5[0,0] = 45042.52777777778
$row5[0,1] = 4.71
$row5[0,2] = 3.47
$row5[0,3] = 0.72
$row5[0,4] = 10.35
$row5[0,5] = 8.01
$row5[0,6] = 3.97
$row5[0,7] = 12.62
$row5[0,8] = 5.82
$row5[0,9] = 2.32
$row5[0,10] = 3.71
$row5[0,11] = 4.08
$row5[0,12] = 4.15
$row5[0,13] = 1.11
$row5[0,14] = 3.76
$row5[0,15] = 5.15
$row5[0,16] = 3.35
$row5[0,17] = 0.71
$row5[0,18] = 0.35
$row5[0,19] = 50.28
$row5[0,20] = 10.3
$row5[0,21] = 3.47
$row5[0,22] = 6.62
$row5[0,23] = 3.51
$row5[0,24] = 0.8100000000000001
$row5[0,25] = 6
$row5[0,26] = 3.07
$row5[0,27] = 2.79
$row5[0,28] = 3.27
$row5[0,29] = 4.32
$row5[0,30] = 0.55
$row5[0,31] = 11.18
$row5[0,32] = 1.85
$row5[0,33] = 4.34
$ws.Range("A5:AH5").Value2 = $row5

# --- Remove the now-obsolete last data row (row 6); used range shrinks to A1:AH5 ---
$ws.Rows.Item(6).Delete() | Out-Null

# --- Resize columns whose content now needs a different width.
#     COM ColumnWidth is in "characters"; the stored OOXML width is ColumnWidth + 5/6,
#     so subtract 5/6 to land on the exact target integer width. ---
$ws.Columns.Item(2).ColumnWidth = 8 - (5/6)  # column B -> width 8
$ws.Columns.Item(3).ColumnWidth = 8 - (5/6)  # column C -> width 8
$ws.Columns.Item(6).ColumnWidth = 7 - (5/6)  # column F -> width 7
$ws.Columns.Item(7).ColumnWidth = 8 - (5/6)  # column G -> width 8
$ws.Columns.Item(11).ColumnWidth = 8 - (5/6)  # column K -> width 8
$ws.Columns.Item(12).ColumnWidth = 8 - (5/6)  # column L -> width 8
$ws.Columns.Item(13).ColumnWidth = 8 - (5/6)  # column M -> width 8
$ws.Columns.Item(15).ColumnWidth = 8 - (5/6)  # column O -> width 8
$ws.Columns.Item(16).ColumnWidth = 8 - (5/6)  # column P -> width 8
$ws.Columns.Item(17).ColumnWidth = 8 - (5/6)  # column Q -> width 8
$ws.Columns.Item(20).ColumnWidth = 9 - (5/6)  # column T -> width 9
$ws.Columns.Item(22).ColumnWidth = 8 - (5/6)  # column V -> width 8
$ws.Columns.Item(24).ColumnWidth = 8 - (5/6)  # column X -> width 8
$ws.Columns.Item(27).ColumnWidth = 8 - (5/6)  # column AA -> width 8
$ws.Columns.Item(29).ColumnWidth = 8 - (5/6)  # column AC -> width 8
$ws.Columns.Item(30).ColumnWidth = 8 - (5/6)  # column AD -> width 8
